$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-48: update Price (D) and Volume(1h) (E) columns ---
# NumberFormat "@" (Text) is applied to column D before assigning so that
# decimal-looking strings (e.g. "0.07000", "101.10") keep their exact
# text representation instead of being auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.266.05"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.68"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7399"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.67"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3145"
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.17"
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07000"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7810"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07968"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.927.57"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.301"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.06"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.34"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.168.05"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.849"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.22"
$ws.Range("E19").Value = "  -3.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007846"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.140.93"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.673"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.512"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.15"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.01"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1273"
$ws.Range("E28").Value = "  -6.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.122"
$ws.Range("E29").Value = "  -8.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.357"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.547"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.332"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.087"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05216"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.305"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7514"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.760"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01950"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.799"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.397"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.88"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4508"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.945"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.798"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8325"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.938"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.10"
$ws.Range("E48").Value = "  -1.51%  "

# --- Rows 49-51: RocketPoolETH row removed, remaining rows shift up, Maker appended ---
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.40"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1217"
$ws.Range("E50").Value = "  +2.30%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "942.76"
$ws.Range("E51").Value = "  -4.99%  "
